# Adds new class/suspension dates to the "Days with Classes" column of the
# money-tracking table (2017_4) for five students, and relocates the
# Word-managed "_GoBack" bookmark to the last-edited cell (Eduardo Farias'),
# matching what Word itself does when the most recent edit happens there.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellParagraphXml([int]$row, [int]$col, [string]$innerBodyXml) {
    $cellRange = $t.Cell($row, $col).Range
    $pkg = '<?xml version="1.0"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           $innerBodyXml +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    [void]$cellRange.InsertXML($pkg)
}

# Row 2: Eduardo Farias -> append ", 17/4, 10/4" and move the "_GoBack"
# bookmark here (it was previously on the Marcela Becerra row).
$eduardoXml = '<w:body>' +
    '<w:p w:rsidR="008E4C32" w:rsidRDefault="003E15B3">' +
    '<w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' +
    '<w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&#190;, 6/4</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, 17/4, 10/4</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p></w:body>'
Set-CellParagraphXml 2 2 $eduardoXml

# Row 3: Jorge Marin -> append ", 17/4, 21/4"
$jorgeXml = '<w:body>' +
    '<w:p w:rsidR="007A0B59" w:rsidRDefault="003E15B3">' +
    '<w:pPr><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' +
    '<w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&#190;, 7/4</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, 17/4, 21/4</w:t></w:r>' +
    '</w:p></w:body>'
Set-CellParagraphXml 3 2 $jorgeXml

# Row 4: Marcela Becerra -> append ", 21/4" and drop the "_GoBack" bookmark
# that used to sit right after ", 12/4" (it moved to the Eduardo row above).
$marcelaXml = '<w:body>' +
    '<w:p w:rsidR="004D6817" w:rsidRDefault="003D0E23">' +
    '<w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' +
    '<w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>5/4</w:t></w:r>' +
    '<w:r w:rsidR="003E15B3"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, 7/4</w:t></w:r>' +
    '<w:r w:rsidR="00933002"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, 12/4</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, 21/4</w:t></w:r>' +
    '</w:p></w:body>'
Set-CellParagraphXml 4 2 $marcelaXml

# Row 5: Sergio Osorio -> append ", 21/4 (1.5hours)"
$sergioXml = '<w:body>' +
    '<w:p w:rsidR="003E15B3" w:rsidRDefault="003E15B3">' +
    '<w:pPr><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' +
    '<w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>7/4(3hours)</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, 21/4 (1.5hours)</w:t></w:r>' +
    '</w:p></w:body>'
Set-CellParagraphXml 5 2 $sergioXml

# Row 6: Cristian Neira -> append ", 21/4"
$cristianXml = '<w:body>' +
    '<w:p w:rsidR="003E15B3" w:rsidRDefault="003E15B3">' +
    '<w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' +
    '<w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>7/4</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, 21/4</w:t></w:r>' +
    '</w:p></w:body>'
Set-CellParagraphXml 6 2 $cristianXml

Write-Output "Applied all five cell updates."
